$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 550
$ws.Range("B3").Value = 300
$ws.Range("B4").Value = 500
$ws.Range("B5").Value = 70
$ws.Range("B6").Value = 90
$ws.Range("B8").Value = 250
